$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dtFmt = "yyyy-mm-dd h:mm:ss"

# Update existing row 196 B/C values with updated precision (timestamps recomputed on re-extraction)
$ws.Range("B196").Value = 45030.43379540509
$ws.Range("C196").Value = 45030.4356243287

# Add new row 197 - full task record
$ws.Range("A197").Value = "Sachin.J"
$ws.Range("B197").Value = 45030.48791054398
$ws.Range("B197").NumberFormat = $dtFmt
$ws.Range("C197").Value = 45030.48954752315
$ws.Range("C197").NumberFormat = $dtFmt
$ws.Range("D197").Value = "Task Completed"

# Add new row 198 (only start time recorded so far)
$ws.Range("A198").Value = "Sachin.J"
$ws.Range("B198").Value = 45030.76237565972
$ws.Range("B198").NumberFormat = $dtFmt

# Add new row 199 (only start time recorded so far)
$ws.Range("A199").Value = "Sachin.J"
$ws.Range("B199").Value = 45030.76890122537
$ws.Range("B199").NumberFormat = $dtFmt
